# Generate Report for Handback
# Updates the "Correspond Handback DateTime" (and, where a new handoff round
# occurred, the "Correspond Handoff File" / "Correspond Handoff Datetime" /
# "Correspond Handback File" columns) on the per-locale report sheets.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row for b07ffab4-... only gets new handoff/handback datetimes ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 14:38:19"
$wsZhCn.Range("H3").Value = "2016-03-18 14:38:35"

# --- de-de sheet: both rows get a new handoff round recorded ---
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2 (4b198396-... file)
$wsDeDe.Range("D2").Value = "4b198396-4f25-4d29-bd36-dbe4c91b5860.022a54277e28870b9c08acfbc58491b14b07b471.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-18 14:37:28"
$wsDeDe.Range("G2").Value = "4b198396-4f25-4d29-bd36-dbe4c91b5860.022a54277e28870b9c08acfbc58491b14b07b471.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-03-18 14:37:54"

# Row 3 (b07ffab4-... file)
$wsDeDe.Range("D3").Value = "b07ffab4-5540-460d-9686-9f583923cf1a.30512969074e4414971dd16c36f2c3c0adb067ab.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-18 14:38:22"
$wsDeDe.Range("G3").Value = "b07ffab4-5540-460d-9686-9f583923cf1a.30512969074e4414971dd16c36f2c3c0adb067ab.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-03-18 14:38:40"
